$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 460
$wsExhibit.Range("F3").Value = 5537
$wsExhibit.Range("F5").Value = 64
$wsExhibit.Range("F6").Value = 82
$wsExhibit.Range("F10").Value = 15

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 460
$wsAll.Range("F3").Value = 5537
$wsAll.Range("F6").Value = 64
$wsAll.Range("F7").Value = 82
$wsAll.Range("F12").Value = 15
